$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update input assumptions for the Black-Scholes calculator
$ws.Range("G10").Value = 135.4
$ws.Range("G11").Value = 140
$ws.Range("G12").Formula = "=4/365"
$ws.Range("G14").Value = 0.4622

# Update the view state to match the saved selection/scroll position
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("G18").Select()
